$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 71.7768669128418, 96252),
    @(3, 71.72727584838867, 104700),
    @(4, 72.11661338806152, 105212),
    @(5, 73.54855537414551, 105468),
    @(6, 123.4285831451416, 105596),
    @(7, 130.6262016296387, 105596),
    @(8, 71.7930793762207, 105596),
    @(9, 74.08308982849121, 105596),
    @(10, 71.31147384643555, 105596),
    @(11, 92.52834320068359, 105596),
    @(12, 15893.02802085876, 95864),
    @(13, 11923.79331588745, 105208),
    @(14, 11435.37735939026, 105720),
    @(15, 11683.5036277771, 105848),
    @(16, 13989.67957496643, 105848),
    @(17, 15115.25297164917, 105848),
    @(18, 11457.23628997803, 105848),
    @(19, 11508.92782211304, 105848),
    @(20, 11701.92646980286, 105848),
    @(21, 11540.26675224304, 105848)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
